$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 7.424215059809214
$ws.Range("E2").Value = 14.69926045795804

$ws.Range("C3").Value = -8.992252553594248
$ws.Range("E3").Value = -19.76480035196674

$ws.Range("C4").Value = 7.007132997505194
$ws.Range("E4").Value = 26.6762552377682

$ws.Range("C5").Value = 8.86644397614711
$ws.Range("E5").Value = 5.982869268853186

$ws.Range("C6").Value = 3.001306214623578
$ws.Range("E6").Value = 0.08911682035466217

$ws.Range("C7").Value = -2.904769335987201
$ws.Range("E7").Value = -8.13908495190001

$ws.Range("C8").Value = 6.148460028297587
$ws.Range("E8").Value = 9.131012060398724

$ws.Range("C9").Value = 5.171596082708629
$ws.Range("E9").Value = 5.986046065677453

$ws.Range("C10").Value = 4.337699953939178
$ws.Range("E10").Value = 5.703761500281579

$ws.Range("C11").Value = 4.073887526082065
$ws.Range("E11").Value = 1.602827009554897

$ws.Range("C12").Value = 2.281540236993274
$ws.Range("E12").Value = -4.308894244053674

$ws.Range("C13").Value = 3.436682959168125
$ws.Range("E13").Value = 2.82953744009995

$ws.Range("C14").Value = -2.305533699949835
$ws.Range("E14").Value = -6.248031846400004

$ws.Range("C15").Value = 5.121832664816339
$ws.Range("E15").Value = 12.25634856540583

$ws.Range("C16").Value = 8.333485306093348
$ws.Range("E16").Value = 7.042341419899389

$ws.Range("C17").Value = 0.08405665459807476
$ws.Range("E17").Value = 0.4484505192704713

$ws.Range("C18").Value = -2.055826728150212
$ws.Range("E18").Value = 3.376972582720295

$ws.Range("C19").Value = 1.378024997308636
$ws.Range("E19").Value = 0.4446401485209472
